$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the previously used range to make room for the new layout
$ws.Cells.Clear()

# Row 1: header
$ws.Range("A1").Value = "with real capacities"
$ws.Range("F1").Value = "unit"

# Row 2: column headers
$ws.Range("B2").Value = "Benchmark"
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 2
$ws.Range("E2").Value = 3

# Row 3: Huber Loss
$ws.Range("B3").Value = "Huber Loss"
$ws.Range("C3").Value = 27.05
$ws.Range("D3").Value = 25.21
$ws.Range("E3").Value = 27.65

# Row 4: Mean Absolute Error (MAE)
$ws.Range("B4").Value = "Mean Absolute Error (MAE)"
$ws.Range("C4").Value = 27.52
$ws.Range("D4").Value = 25.67
$ws.Range("E4").Value = 28.11
$ws.Range("F4").Value = "MW"

# Row 5: Mean Squared Error (MSE)
$ws.Range("B5").Value = "Mean Squared Error (MSE)"
$ws.Range("C5").Value = 2941.27
$ws.Range("D5").Value = 2681.88
$ws.Range("E5").Value = 3077
$ws.Range("F5").Value = "MW²"

# Row 6: Root Mean Squared Error (RMSE)
$ws.Range("B6").Value = "Root Mean Squared Error (RMSE)"
$ws.Range("C6").Value = 54.23
$ws.Range("D6").Value = 51.79
$ws.Range("E6").Value = 55.47

# Row 7: normalised header
$ws.Range("A7").Value = "normalised"

# Row 8: Huber Loss
$ws.Range("B8").Value = "Huber Loss"
$ws.Range("C8").Value = 0.0309
$ws.Range("D8").Value = 0.0285
$ws.Range("E8").Value = 0.0309

# Row 9: Mean Absolute Error (MAE)
$ws.Range("B9").Value = "Mean Absolute Error (MAE)"
$ws.Range("C9").Value = 0.1607
$ws.Range("D9").Value = 0.1522
$ws.Range("E9").Value = 0.1623

# Row 10: Mean Squared Error (MSE)
$ws.Range("B10").Value = "Mean Squared Error (MSE)"
$ws.Range("C10").Value = 0.0617
$ws.Range("D10").Value = 0.0571
$ws.Range("E10").Value = 0.0619

# Row 11: Root Mean Squared Error (RMSE)
$ws.Range("B11").Value = "Root Mean Squared Error (RMSE)"
$ws.Range("C11").Value = 0.2484
$ws.Range("D11").Value = 0.239
$ws.Range("E11").Value = 0.2488

# Row 12: Model header
$ws.Range("B12").Value = "Model (100.000 data points, with cut-off of 0.01 values)"

# Row 13: Huber Loss
$ws.Range("B13").Value = "Huber Loss"
$ws.Range("C13").Value = 0.0143

# Row 14: Mean Absolute Error (MAE)
$ws.Range("B14").Value = "Mean Absolute Error (MAE)"
$ws.Range("C14").Value = 0.1212

# Row 15: Mean Squared Error (MSE)
$ws.Range("B15").Value = "Mean Squared Error (MSE)"
$ws.Range("C15").Value = 0.0287

# Row 16: Root Mean Squared Error (RMSE)
$ws.Range("B16").Value = "Root Mean Squared Error (RMSE)"
$ws.Range("C16").Value = 0.1694

# Restore the selection as in the final workbook
$ws.Range("E15").Select()
